# Generate Report for Handoff
# A new handoff round was generated: the source file was renamed from the
# "9a9450c5-..." GUID to "a195be73-...", a fresh xliff handoff was produced
# (new xliff hash, new handoff timestamps), and the previous handback
# file/date fields are reset since no handback has happened yet for this
# new round.

$wb = $excel.ActiveWorkbook

$newGuid = "a195be73-b240-4d7b-b149-635339d546cd"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")

$ovw.Range("A2").Value = "$newGuid.md"
$ovw.Range("B2").Value = "e2e\$newGuid.md"
$ovw.Range("G2").Value = "2016-09-02 05:04:45"

foreach ($hl in $ovw.Hyperlinks) {
    if ($hl.Range.Address() -eq '$B$2') {
        $hl.TextToDisplay = "e2e\$newGuid.md"
    }
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("A2").Value = "$newGuid.md"
$zhcn.Range("G2").Value = "$newGuid.39eb592dc5a3699096ab914087d9ccfd7650e580.zh-cn.xlf"
$zhcn.Range("H2").Value = "2016-09-02 05:04:36"

foreach ($hl in $zhcn.Hyperlinks) {
    if ($hl.Range.Address() -eq '$A$2') {
        $hl.TextToDisplay = "$newGuid.md"
    }
    if ($hl.Range.Address() -eq '$I$2') {
        $hl.Delete()
    }
}

$zhcn.Range("I2").Style = "Normal"
$zhcn.Range("I2").Value = ""
$zhcn.Range("J2").Value = ""
$zhcn.Range("K2").Value = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("A2").Value = "$newGuid.md"
$dede.Range("G2").Value = "$newGuid.39eb592dc5a3699096ab914087d9ccfd7650e580.de-de.xlf"
$dede.Range("H2").Value = "2016-09-02 05:04:45"

foreach ($hl in $dede.Hyperlinks) {
    if ($hl.Range.Address() -eq '$A$2') {
        $hl.TextToDisplay = "$newGuid.md"
    }
    if ($hl.Range.Address() -eq '$I$2') {
        $hl.Delete()
    }
}

$dede.Range("I2").Style = "Normal"
$dede.Range("I2").Value = ""
$dede.Range("J2").Value = ""
$dede.Range("K2").Value = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Column width adjustments (zh-cn / de-de columns I & J) - these columns
# no longer need to show full file names (the Latest Target / Latest
# Handback columns are now blank for a fresh handoff), so they are
# narrowed down from the default 40 characters.
# ---------------------------------------------------------------------
foreach ($ws in @($zhcn, $dede)) {
    $ws.Columns.Item(9).ColumnWidth = 17.83
    $ws.Columns.Item(10).ColumnWidth = 20.83
}
